# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1580
    3  = 53
    4  = 1039
    7  = 2743
    9  = 1759
    10 = 186
    11 = 78
    12 = 609
    13 = 30
    14 = 21
    15 = 148
    16 = 81
    17 = 85
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
